# Fruta / hortaliza, semanal
# Insert a new data row above the current row 16 (shifting existing rows
# 16-39 down to 17-40), then populate the new row 16 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16..39 down to 17..40 by inserting a new row at 16.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16.
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value = "Maule"
$ws.Cells.Item(16, 4).Value = 44645
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = "Frutos de pepita"
$ws.Cells.Item(16, 9).Value = 100104003
$ws.Cells.Item(16, 10).Value = "Membrillo"
$ws.Cells.Item(16, 11).Value = "Champion"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 10000
$ws.Cells.Item(16, 15).Value = 10000
$ws.Cells.Item(16, 16).Value = 10000
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 556
$ws.Cells.Item(16, 20).Value = 18

# Make sure the D column keeps the date number format used by the rest of
# the sheet's date column (style copied down automatically by the row
# insert, but set explicitly to be safe).
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
